$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 85, shifting existing rows 85:90 down to 86:91
$ws.Rows.Item(85).Insert()

# Populate the new row 85 with the same "constant" columns as the rest of
# the Berenjena / Vega Monumental Concepcion block, plus the new record's
# specific values.
$ws.Cells.Item(85, 1).Value = 11
$ws.Cells.Item(85, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(85, 3).Value = "Bíobío"
$ws.Cells.Item(85, 4).Value = 44714
$ws.Cells.Item(85, 5).Value = 8
$ws.Cells.Item(85, 6).Value = 100112001
$ws.Cells.Item(85, 7).Value = "Berenjena"
$ws.Cells.Item(85, 8).Value = "Sin especificar"
$ws.Cells.Item(85, 9).Value = "Primera"
$ws.Cells.Item(85, 10).Value = 270
$ws.Cells.Item(85, 11).Value = 5500
$ws.Cells.Item(85, 12).Value = 6000
$ws.Cells.Item(85, 13).Value = 5722
$ws.Cells.Item(85, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(85, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(85, 16).Value = 95
$ws.Cells.Item(85, 17).Value = 60
$ws.Cells.Item(85, 18).Value = "Hortaliza"
